$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: update total_venda (B2)
$ws.Range("B2").Value = 15517.98

# Row 3: update Dia (A3) and total_venda (B3)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 27866.61

# Row 20: update total_venda (B20)
$ws.Range("B20").Value = 15253.9

# Row 24: update total_venda (B24)
$ws.Range("B24").Value = 114289.79
